$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "63.796.09" with dot-
# grouped thousands, "5.40"/"164.40" with a significant trailing zero,
# "0.0000164" that Excel would otherwise re-render in scientific notation).
# They must stay plain text exactly as authored, so force each touched D cell
# to Text format before writing, then drop back to the default "Normal" style
# so no stray number-format/style is left on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.796.09'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.621.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  +4.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.394'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.33%  '
$ws.Range('E11').Value = '  +3.56%  '
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.091.30'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.687.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000164'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +10.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.638.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '348.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +1.84%  '
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.06%  '
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  +3.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '552.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.162'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0889'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.27%  '
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.40'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.17'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.05%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.417'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.69%  '
$ws.Range('E38').Value = '  +2.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '167.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('E43').Value = '  +4.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0585'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.636'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0968'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.19'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('E51').Value = '  +18.93%  '
